$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row100 = @("Globo", "RJ TV 2", "Obras", "2025-04-08T18:47", "Negativo", "no meio do caminho tinha um buraco...Moradores do Parque Esplanada e de Guarus reclamam das crateras nas vias. Depoimento de morador, alegando que acidentes estão acontecendo no local. Acidentes já aconteceram e os buracos estão atraindo ratos e cobras. Em Guarus, buraco próximo à delegacia da região. Depoimento de morador do bairro. *com nota da prefeitura*")

$row101 = @("Globo", "RJ TV 2", "Saúde", "2025-04-09T18:49", "Neutro", "Profissionais de Campos temem demissão. Destino incerto de RPAs da atenção domiciliar preocupa familiares de pacientes. Repórter *ao vivo*. Profissionais se reuniram de manhã em frente à Prefeitura de Campos. Com Reforma Administrativa, muitos RPAs serão demitidos. Alguns já foram. Reforma administrativa ainda não foi votada por decisão judicial. Prefeitura já recorreu. Entrevista com profissionais que participaram da manifestação. *com nota da prefeitura*. Secretaria de Saúde informou que o SAD está mantido e não há interrupção de atendimento. ")

for ($i = 0; $i -lt $row100.Length; $i++) {
    $ws.Cells.Item(100, $i + 1).Value = $row100[$i]
}

for ($i = 0; $i -lt $row101.Length; $i++) {
    $ws.Cells.Item(101, $i + 1).Value = $row101[$i]
}
